$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 11: Escherichia coli CFT073 (UPEC) ----
$ws.Rows.Item(11).RowHeight = 13.5
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("A11").Value = 199310.4
$ws.Range("B11").Value = "CFT073"
$ws.Range("C11").Value = "Escherichia coli Cft073"
$ws.Range("D11").Value = "InPec: UPEC"
$ws.Range("E11").Value = "AE014075"

# ---- Row 12: Escherichia coli UMN026 (UPEC) ----
$ws.Rows.Item(12).RowHeight = 13.5
$ws.Range("D10").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B12").PasteSpecial(-4122)

$ws.Range("A12").Value = 585056.7
$ws.Range("C12").Value = "Escherichia coli UMN026"
$ws.Range("C12").Font.Bold = $true
$ws.Range("B12").Value = "UMN026"
$ws.Range("D12").Value = "InPec: UPEC"
$ws.Range("E12").Value = "CU928163"

# ---- Row 13: Shigella flexneri 2a str. 2457T (EHEC) ----
$ws.Rows.Item(13).RowHeight = 13.5
$ws.Range("D10").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)

$ws.Range("A13").Value = 198215.6
$ws.Range("E13").Value = "AE014073"
$ws.Range("C13").Value = "Shigella flexneri 2a str. 2457T"
$ws.Range("C13").Font.Bold = $true
$ws.Range("B13").Value = "Shigella 2457T"
$ws.Range("D13").Value = "InPec: EHEC"

# ---- Row 14: Escherichia coli SMS-3-5 ----
$ws.Rows.Item(14).RowHeight = 13.5
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("A14").Value = 439855.1
$ws.Range("E14").Value = "CP000970"
$ws.Range("C14").Value = "Escherichia coli SMS-3-5"
$ws.Range("C14").Font.Bold = $true
$ws.Range("B14").Value = "SMS-3-5"

# ---- Row 15: Escherichia coli ATCC 8739 ----
$ws.Rows.Item(15).RowHeight = 13.5
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("A15").Value = 481805.6
$ws.Range("E15").Value = "CP000946"
$ws.Range("C15").Value = "Escherichia coli ATCC 8739"
$ws.Range("C15").Font.Bold = $true
$ws.Range("B15").Value = "ATCC 8739"

# ---- Row 16: stray touched cell ----
$ws.Rows.Item(16).RowHeight = 13.5
$ws.Range("B16").Value = "x"
$ws.Range("B16").ClearContents()
$ws.Range("B16").Style = "Normal"

# ---- Final selection ----
$null = $ws.Range("C14").Select()
